# Append the 10/22/2025 profit-allocation row (A51:C51) to the sheet,
# matching the format of the existing data rows (plain text date,
# numeric BTC/KAS shares, default "Normal" cell style).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 51

# Force the date cell to be stored as literal text (not auto-converted
# to a date serial number) by applying a text number format before
# assigning the value, then restore the default "Normal" style so the
# new row's formatting matches the other data rows (no explicit style).
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "10/22/2025"
$ws.Cells.Item($row, 1).Style = "Normal"

$ws.Cells.Item($row, 2).Value = 0.1958319825464067
$ws.Cells.Item($row, 3).Value = 0.8041680174535933
